$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 3 (the totals row),
# shifting it down to row 5.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Row 2: update the first data row in place
$ws.Range("A2").Value = "901/LF/FES "
$ws.Range("B2").Value = "Logement de fonction"
$ws.Range("C2").Value = "BJ36877"
$ws.Range("D2").Value = "CHARIJI ABDELLAH"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = "--"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "--"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 1200
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 10800

# Row 3: new data row
$ws.Range("A3").Value = "901/LF/FES "
$ws.Range("B3").Value = "Logement de fonction"
$ws.Range("C3").Value = "BJ36877"
$ws.Range("D3").Value = "CHARIJI ABDELLAH"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 6000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 11400

# Row 4: new data row
$ws.Range("A4").Value = "901/FES "
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "J207703"
$ws.Range("D4").Value = "ACHENGLI LAILA"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 9500

# Row 5: totals row (A-G stay blank/space, H-O updated totals)
$ws.Range("H5").Value = 11000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1100
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 1200
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 31700
